$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (raw serial write, keep date format)
$ws.Range("A1").Value2 = 45309

# Step 2: update prices in column D for the two price tables
$ws.Range("D10").Value2 = 3231
$ws.Range("D11").Value2 = 3913
$ws.Range("D12").Value2 = 3640
$ws.Range("D13").Value2 = 4459
$ws.Range("D14").Value2 = 5187
$ws.Range("D15").Value2 = 3822
$ws.Range("D16").Value2 = 6753
$ws.Range("D17").Value2 = 324

$ws.Range("D25").Value2 = 4050
$ws.Range("D26").Value2 = 7080
$ws.Range("D27").Value2 = 5460
$ws.Range("D28").Value2 = 7470
$ws.Range("D29").Value2 = 5100
$ws.Range("D30").Value2 = 8477
$ws.Range("D31").Value2 = 6734
$ws.Range("D32").Value2 = 8750
